$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
  @(9,9),
  @(9,9),
  @(9,9),
  @(9,9),
  @(9,9),
  @(9,9),
  @(7,8),
  @(7,7),
  @(7,7),
  @(7,8),
  @(7,7),
  @(6,7),
  @(8,8),
  @(9,9),
  @(9,9),
  @(6,7),
  @(5,6),
  @(8,9),
  @(7,8),
  @(6,6),
  @(8,9),
  @(5,6),
  @(6,6),
  @(5,7),
  @(8,8),
  @(7,7),
  @(5,6),
  @(8,9),
  @(6,7),
  @(8,8),
  @(8,8),
  @(8,9),
  @(7,8),
  @(4,6),
  @(7,7),
  @(11,6),
  @(8,8),
  @(7,7),
  @(2,3),
  @(7,7),
  @(4,4),
  @(4,4),
  @(5,5),
  @(4,4),
  @(5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
